$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.639840960502625
$ws.Range("B1").Value = 3.195481061935425
$ws.Range("C1").Value = 5.078707218170166
$ws.Range("D1").Value = 1.417985081672668
$ws.Range("E1").Value = 0.8301410675048828
